# Applies the "Updated cryptos list" refresh: new Price (D) / Volume(1h) (E)
# values for the rows whose figures moved in this run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.672.49"
$ws.Range("E2").Value = "  +5.64%  "

$ws.Range("D3").Value = "3.633.85"
$ws.Range("E3").Value = "  +5.62%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.81"
$ws.Range("E5").Value = "  +2.26%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "194.79"
$ws.Range("E6").Value = "  +3.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.644"
$ws.Range("E7").Value = "  +2.51%  "

$ws.Range("D8").Value = "3.626.20"
$ws.Range("E8").Value = "  +5.63%  "

$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("E10").Value = "  +7.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.679"
$ws.Range("E11").Value = "  +5.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.01"
$ws.Range("E12").Value = "  +1.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000313"
$ws.Range("E13").Value = "  +13.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.93"
$ws.Range("E14").Value = "  +5.14%  "

$ws.Range("D15").Value = "4.217.47"
$ws.Range("E15").Value = "  +5.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.49"
$ws.Range("E16").Value = "  +8.77%  "

$ws.Range("D17").Value = "3.638.09"
$ws.Range("E17").Value = "  +5.79%  "

$ws.Range("D18").Value = "70.687.62"
$ws.Range("E18").Value = "  +5.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.72"
$ws.Range("E19").Value = "  +5.64%  "

$ws.Range("E20").Value = "  +2.71%  "

$ws.Range("E21").Value = "  +4.23%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "488.53"
$ws.Range("E22").Value = "  +2.65%  "

$ws.Range("E23").Value = "  +13.56%  "

$ws.Range("E24").Value = "  -3.10%  "

$ws.Range("E25").Value = "  +3.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.34"
$ws.Range("E26").Value = "  +2.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.17"
$ws.Range("E27").Value = "  +6.85%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  +4.22%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.57"
$ws.Range("E29").Value = "  +6.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.95"
$ws.Range("E30").Value = "  +8.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.80"
$ws.Range("E31").Value = "  +5.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.122"
$ws.Range("E32").Value = "  +9.75%  "

$ws.Range("E33").Value = "  +4.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.33"
$ws.Range("E34").Value = "  +2.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "611.41"
$ws.Range("E35").Value = "  +2.33%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.54"
$ws.Range("E36").Value = "  +9.65%  "

$ws.Range("D37").Value = "0.0₃0839"
$ws.Range("E37").Value = "  +11.85%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.411"
$ws.Range("E38").Value = "  +5.55%  "

$ws.Range("E39").Value = "  +1.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.12%  "

$ws.Range("D42").Value = "3.322.58"
$ws.Range("E42").Value = "  +3.97%  "

$ws.Range("E43").Value = "  +17.47%  "

$ws.Range("E44").Value = "  +9.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.84"
$ws.Range("E45").Value = "  +9.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0459"
$ws.Range("E46").Value = "  +6.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.68"
$ws.Range("E47").Value = "  +12.79%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.33"
$ws.Range("E48").Value = "  +2.53%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.139"
$ws.Range("E49").Value = "  +3.43%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.24"
$ws.Range("E50").Value = "  +1.72%  "

$ws.Range("E51").Value = "  +0.09%  "
